# repull data, push all data, mean calculation
# Update the dSF (column F) values on Sheet1 to reflect the re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    5  = 0
    6  = -1
    10 = -3
    12 = -3
    19 = -2
    23 = 1
    24 = 3
    29 = -1
    30 = -2
    32 = 0
    33 = 1
    36 = 1
    37 = 0
    38 = -1
    42 = 2
    44 = -7
    48 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
